# Insert a new data row for 2026/01/22 (16:00 slot) right before the
# 2026/12/29 block, pushing the existing rows 693..734 down to 694..735.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 693 downward by inserting a new blank row at 693.
$ws.Rows.Item(693).Insert()

# Populate the newly inserted row with the continuation-of-day data.
# The date/day-of-week columns hold plain text that merely looks like a
# date (matching the rest of the sheet), so force the cell to Text before
# assignment to stop Excel's automatic date recognition, then drop the
# now-unneeded number format so the cell keeps the sheet's default style.
$ws.Cells.Item(693, 1).NumberFormat = "@"
$ws.Cells.Item(693, 1).Value = "2026/01/22"
$ws.Cells.Item(693, 1).ClearFormats()

$ws.Cells.Item(693, 2).Value = "木"
$ws.Cells.Item(693, 3).Value = 16
$ws.Cells.Item(693, 4).Value = 201
